$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Ngày bắt đầu" (C) and "Ngày kết thúc" (D) columns entirely.
# This shifts "Tổng số môn học" (old E) into C and "Tổng số tín chỉ" (old F) into D.
$ws.Range("C:D").Delete()

# Convert the "Học kỳ" values from text ("HK1"/"HK2") to plain numbers (1/2).
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
